$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.2109803333333333
$ws.Range("H2").Value = 0.632941
$ws.Range("I2").Value = 0.2372796149188938
$ws.Range("J2").Value = 0.2372796149188938
$ws.Range("M2").Value = 0.00535
$ws.Range("N2").Value = 0.01605
$ws.Range("O2").Value = 0.003591913026022235
$ws.Range("P2").Value = 0.003591913026022235
$ws.Range("Q2").Value = 0.001128744783333333
$ws.Range("R2").Value = 0.01015870305
$ws.Range("S2").Value = 0.0008522877396367145
$ws.Range("T2").Value = 0.0008522877396367145
$ws.Range("G3").Value = 0.2109803333333333
$ws.Range("H3").Value = 0.632941
$ws.Range("I3").Value = 0.2372796149188938
$ws.Range("J3").Value = 0.2372796149188938
$ws.Range("O3").Value = 0.9964080869739778
$ws.Range("P3").Value = 0.9964080869739778
$ws.Range("Q3").Value = 0.3131173895623333
$ws.Range("R3").Value = 2.818056506061
$ws.Range("S3").Value = 0.2364273271792571
$ws.Range("T3").Value = 0.2364273271792571
$ws.Range("I4").Value = 0.4723879002358022
$ws.Range("J4").Value = 0.4723879002358022
$ws.Range("M4").Value = 0.00535
$ws.Range("N4").Value = 0.01605
$ws.Range("O4").Value = 0.003591913026022235
$ws.Range("P4").Value = 0.003591913026022235
$ws.Range("Q4").Value = 0.0022471605
$ws.Range("R4").Value = 0.0202244445
$ws.Range("S4").Value = 0.00169677625219227
$ws.Range("T4").Value = 0.00169677625219227
$ws.Range("I5").Value = 0.4723879002358022
$ws.Range("J5").Value = 0.4723879002358022
$ws.Range("O5").Value = 0.9964080869739778
$ws.Range("P5").Value = 0.9964080869739778
$ws.Range("Q5").Value = 0.6233694632099999
$ws.Range("R5").Value = 5.610325168889999
$ws.Range("S5").Value = 0.4706911239836099
$ws.Range("T5").Value = 0.4706911239836099
$ws.Range("I6").Value = 0.2903324848453041
$ws.Range("J6").Value = 0.2903324848453041
$ws.Range("M6").Value = 0.00535
$ws.Range("N6").Value = 0.01605
$ws.Range("O6").Value = 0.003591913026022235
$ws.Range("P6").Value = 0.003591913026022235
$ws.Range("Q6").Value = 0.00138111855
$ws.Range("R6").Value = 0.01243006695
$ws.Range("S6").Value = 0.001042849034193251
$ws.Range("T6").Value = 0.001042849034193251
$ws.Range("I7").Value = 0.2903324848453041
$ws.Range("J7").Value = 0.2903324848453041
$ws.Range("O7").Value = 0.9964080869739778
$ws.Range("P7").Value = 0.9964080869739778
$ws.Range("Q7").Value = 0.383126674371
$ws.Range("S7").Value = 0.2892896358111108
$ws.Range("T7").Value = 0.2892896358111108
